$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1244.625
$ws.Range("I15").Value = 1244.625
$ws.Range("K15").Value = 3733.875
$ws.Range("M15").Value = -3564.875
$ws.Range("H19").Value = 2352.125
$ws.Range("I19").Value = 2287.1538
$ws.Range("J19").Value = 2633.6667
$ws.Range("K19").Value = 2287.1538
$ws.Range("L19").Value = 2633.6667
$ws.Range("M19").Value = -2112.1538
$ws.Range("N19").Value = -2983.6667
$ws.Range("H129").Value = 2786.111
$ws.Range("I129").Value = 1030.1666
$ws.Range("J129").Value = 6298
$ws.Range("K129").Value = 3090.4998
$ws.Range("L129").Value = 18894
$ws.Range("M129").Value = 1909.5002
$ws.Range("N129").Value = -28894
$ws.Range("H132").Value = 3879.6758
$ws.Range("I132").Value = 1388.9354
$ws.Range("K132").Value = 4166.8062
$ws.Range("M132").Value = -1636.8062
$ws.Range("H137").Value = 2819.8386
$ws.Range("I137").Value = 2236.6316
$ws.Range("J137").Value = 3743.25
$ws.Range("K137").Value = 6709.8948
$ws.Range("L137").Value = 11229.75
$ws.Range("M137").Value = -4159.8948
$ws.Range("N137").Value = -16329.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14658.113
$ws.Range("I32").Value = 13610.614
$ws.Range("K32").Value = 13610.614
$ws.Range("M32").Value = -13323.614
$ws.Range("H45").Value = 2693.5454
$ws.Range("I45").Value = 1471.8334
$ws.Range("K45").Value = 1471.8334
$ws.Range("M45").Value = -1094.8334
$ws.Range("H61").Value = 4854.9414
$ws.Range("I61").Value = 3538.1428
$ws.Range("J61").Value = 11000
$ws.Range("K61").Value = 3538.1428
$ws.Range("L61").Value = 11000
$ws.Range("M61").Value = -3326.1428
$ws.Range("N61").Value = -11424
$ws.Range("H74").Value = 19231636
$ws.Range("I74").Value = 22728022
$ws.Range("J74").Value = 1506.5
$ws.Range("K74").Value = 22728022
$ws.Range("L74").Value = 1506.5
$ws.Range("M74").Value = -22727148
$ws.Range("N74").Value = -3254.5
$ws.Range("H77").Value = 19231636
$ws.Range("I77").Value = 22728022
$ws.Range("J77").Value = 1506.5
$ws.Range("K77").Value = 113640110
$ws.Range("L77").Value = 7532.5
$ws.Range("M77").Value = -113635742
$ws.Range("N77").Value = -16268.5
$ws.Range("H102").Value = 5000
$ws.Range("J102").Value = 5000
$ws.Range("L102").Value = 5000
$ws.Range("N102").Value = -8244
$ws.Range("H104").Value = 50741.668
$ws.Range("J104").Value = 50741.668
$ws.Range("L104").Value = 50741.668
$ws.Range("N104").Value = -57729.668
$ws.Range("H122").Value = 5742.5
$ws.Range("I122").Value = 5198.6206
$ws.Range("K122").Value = 15595.8618
$ws.Range("M122").Value = -13145.8618
$ws.Range("H132").Value = 2236.3142
$ws.Range("I132").Value = 2236.3142
$ws.Range("K132").Value = 6708.942599999999
$ws.Range("M132").Value = -4178.942599999999
$ws.Range("H136").Value = 4854.9414
$ws.Range("I136").Value = 3538.1428
$ws.Range("J136").Value = 11000
$ws.Range("K136").Value = 10614.4284
$ws.Range("L136").Value = 33000
$ws.Range("M136").Value = -8064.428400000001
$ws.Range("N136").Value = -38100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2501.625
$ws.Range("I99").Value = 2428.8572
$ws.Range("K99").Value = 2428.8572
$ws.Range("M99").Value = -930.8571999999999
$ws.Range("H103").Value = 9966.333000000001
$ws.Range("J103").Value = 9966.333000000001
$ws.Range("L103").Value = 9966.333000000001
$ws.Range("N103").Value = -12310.333
$ws.Range("H134").Value = 2551.7666
$ws.Range("I134").Value = 1991.0555
$ws.Range("J134").Value = 3392.8333
$ws.Range("K134").Value = 5973.166499999999
$ws.Range("L134").Value = 10178.4999
$ws.Range("M134").Value = -3438.166499999999
$ws.Range("N134").Value = -15248.4999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2153.4
$ws.Range("I22").Value = 2385
$ws.Range("K22").Value = 2385
$ws.Range("M22").Value = -2035
$ws.Range("H52").Value = 46091.668
$ws.Range("J52").Value = 49183.332
$ws.Range("L52").Value = 49183.332
$ws.Range("N52").Value = -49771.332
$ws.Range("H58").Value = 2662.7896
$ws.Range("I58").Value = 2354.0625
$ws.Range("K58").Value = 2354.0625
$ws.Range("M58").Value = -2151.0625
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H104").Value = 59999
$ws.Range("I104").Value = 59999
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 59999
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -57378
$ws.Range("N104").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 75000
$ws.Range("J119").Value = 75000
$ws.Range("L119").Value = 75000
$ws.Range("N119").Value = -84676
$ws.Range("H122").Value = 4190.8423
$ws.Range("I122").Value = 2817.6155
$ws.Range("K122").Value = 8452.8465
$ws.Range("M122").Value = -6002.8465
$ws.Range("H132").Value = 60607610
$ws.Range("I132").Value = 60607610
$ws.Range("K132").Value = 181822830
$ws.Range("M132").Value = -181820300
$ws.Range("H134").Value = 3338.6667
$ws.Range("I134").Value = 2457.0833
$ws.Range("K134").Value = 7371.249899999999
$ws.Range("M134").Value = -4836.249899999999
$ws.Range("H136").Value = 2662.7896
$ws.Range("I136").Value = 2354.0625
$ws.Range("K136").Value = 7062.1875
$ws.Range("M136").Value = -4512.1875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3949.5
$ws.Range("I109").Value = 1290.909
$ws.Range("K109").Value = 3872.727
$ws.Range("M109").Value = -2832.727
$ws.Range("H115").Value = 6603.8335
$ws.Range("I115").Value = 6603.8335
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 19811.5005
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -18636.5005
$ws.Range("N115").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 28664.5
$ws.Range("J95").Value = 28664.5
$ws.Range("L95").Value = 28664.5
$ws.Range("N95").Value = -34156.5
$ws.Range("H96").Value = 39499
$ws.Range("J96").Value = 39499
$ws.Range("L96").Value = 39499
$ws.Range("N96").Value = -44991
$ws.Range("H132").Value = 1893.75
$ws.Range("I132").Value = 1893.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5681.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3151.25
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 964
$ws.Range("I22").Value = 792.55554
$ws.Range("K22").Value = 792.55554
$ws.Range("M22").Value = -497.55554
$ws.Range("H27").Value = 964
$ws.Range("I27").Value = 792.55554
$ws.Range("K27").Value = 792.55554
$ws.Range("M27").Value = -685.55554
$ws.Range("H82").Value = 3785.4
$ws.Range("I82").Value = 2472.8333
$ws.Range("J82").Value = 4660.4443
$ws.Range("K82").Value = 2472.8333
$ws.Range("L82").Value = 4660.4443
$ws.Range("M82").Value = -2111.8333
$ws.Range("N82").Value = -5382.4443
$ws.Range("H85").Value = 3785.4
$ws.Range("I85").Value = 2472.8333
$ws.Range("J85").Value = 4660.4443
$ws.Range("K85").Value = 2472.8333
$ws.Range("L85").Value = 4660.4443
$ws.Range("M85").Value = -1224.8333
$ws.Range("N85").Value = -7156.4443
$ws.Range("H93").Value = 1955.25
$ws.Range("I93").Value = 1968.2667
$ws.Range("J93").Value = 1760
$ws.Range("K93").Value = 1968.2667
$ws.Range("L93").Value = 1760
$ws.Range("M93").Value = -720.2666999999999
$ws.Range("N93").Value = -4256
$ws.Range("H104").Value = 54998
$ws.Range("J104").Value = 54998
$ws.Range("L104").Value = 54998
$ws.Range("N104").Value = -61986
$ws.Range("H132").Value = 2830
$ws.Range("I132").Value = 2796.808
$ws.Range("J132").Value = 2919.7407
$ws.Range("K132").Value = 8390.423999999999
$ws.Range("L132").Value = 8759.222099999999
$ws.Range("M132").Value = -5860.423999999999
$ws.Range("N132").Value = -13819.2221
$ws.Range("H136").Value = 4052.9285
$ws.Range("I136").Value = 2749.0513
$ws.Range("K136").Value = 8247.153900000001
$ws.Range("M136").Value = -5697.153900000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 42328.668
$ws.Range("J101").Value = 42328.668
$ws.Range("L101").Value = 42328.668
$ws.Range("N101").Value = -48818.668
$ws.Range("H136").Value = 4020.157
$ws.Range("I136").Value = 2613.7878
$ws.Range("M136").Value = -5291.3634
